$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''246.88'
$ws.Range("D3").Value = '''22.47'
$ws.Range("D4").Value = '''5.257'
$ws.Range("D5").Value = '''0.05689'
$ws.Range("D6").Value = '''3.416'
$ws.Range("D7").Value = '''6.288'
$ws.Range("D8").Value = '''0.8083'
$ws.Range("D9").Value = '''0.8711'
$ws.Range("D10").Value = '''0.01098'
$ws.Range("D11").Value = '''0.1414'
$ws.Range("D12").Value = '''0.07347'
$ws.Range("D13").Value = '''0.03032'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = '''0.03069'
$ws.Range("E14").Value = '13BitrueCoinBTR'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = '''0.09380'
$ws.Range("E15").Value = '14BitMartTokenBMX'
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").Value = '''3.869'
$ws.Range("E16").Value = '15MCDexMCB'
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").Value = '''0.001585'
$ws.Range("E17").Value = '16BitForexTokenBF'
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").Value = '''0.04766'
$ws.Range("E18").Value = '17CoinExTokenCET'
$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D19").Value = '''0.006386'
$ws.Range("E19").Value = '18TigerCashTCH'
$ws.Range("B20").Value = 'HotbitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D20").Value = '''0.005021'
$ws.Range("E20").Value = '19HotbitTokenHTB'
$ws.Range("B21").Value = 'BitKan'
$ws.Range("C21").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D21").Value = '''0.0009962'
$ws.Range("E21").Value = '20BitKanKAN'
$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D22").Value = '''0.0001500'
$ws.Range("E22").Value = '21NitroExNTX'
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").Value = '''3.691'
$ws.Range("E23").Value = '22LEOLEO'
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").Value = '''2.195'
$ws.Range("E24").Value = '23BTSETokenBTSE'
$ws.Range("B25").Value = 'BitpandaEcosystemToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D25").Value = '''0.3279'
$ws.Range("E25").Value = '24BitpandaEcosystemTokenBEST'
$ws.Range("B26").Value = 'ProBitToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D26").Value = '''0.1341'
$ws.Range("E26").Value = '25ProBitTokenPROB'
$ws.Range("D40").Value = '''0.03919'
$ws.Range("D41").Value = '''0.006799'
$ws.Range("D42").Value = '''0.1065'
$ws.Range("D44").Value = '''0.007486'
$ws.Range("D48").Value = '''0.1950'
